$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column N (2022) data, mirroring the style of column M
$ws.Range("N3").Value = $null
$ws.Range("N4").Value = 2022
$ws.Range("N5").Value = 98.8
$ws.Range("N6").Value = 98
$ws.Range("N7").Value = 96.9

# Copy formatting from column M (the previous last data column) to column N, cell by cell
$ws.Range("M3").Copy() | Out-Null
$ws.Range("N3").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("M4").Copy() | Out-Null
$ws.Range("N4").PasteSpecial(-4122) | Out-Null
$ws.Range("M5").Copy() | Out-Null
$ws.Range("N5").PasteSpecial(-4122) | Out-Null
$ws.Range("M6").Copy() | Out-Null
$ws.Range("N6").PasteSpecial(-4122) | Out-Null
$ws.Range("M7").Copy() | Out-Null
$ws.Range("N7").PasteSpecial(-4122) | Out-Null

# Select cell O4 to match the resulting active selection in the file
$ws.Range("O4").Select() | Out-Null
